$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that currently sits right
#    after "QUESTION2".
# ------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 2. Merge the three separate "China" / "Mandarin" / "Chinese" list
#    bullets into a single paragraph, separated by tab characters,
#    and re-insert the "_GoBack" bookmark right before "Chinese".
# ------------------------------------------------------------------
$china = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "China") {
        $china = $p
        break
    }
}

# Merge "China" paragraph with the following ("Mandarin") paragraph by
# deleting the paragraph mark between them, then put a tab in its place.
$mergePoint1 = $china.Range.End - 1
$d.Range($mergePoint1, $mergePoint1 + 1).Delete()
$d.Range($mergePoint1, $mergePoint1).InsertAfter([char]9)

# Merge the resulting paragraph ("China<tab>Mandarin") with the next one
# ("Chinese") the same way.
$mergePoint2 = $china.Range.End - 1
$d.Range($mergePoint2, $mergePoint2 + 1).Delete()
$d.Range($mergePoint2, $mergePoint2).InsertAfter([char]9)

# Re-add the "_GoBack" bookmark right before "Chinese" (i.e. right after
# the tab that now separates "Mandarin" and "Chinese"). Paragraph.Range
# includes the trailing paragraph mark, so that has to be excluded too.
$bookmarkPos = $china.Range.End - 1 - "Chinese".Length
$d.Bookmarks.Add("_GoBack", $d.Range($bookmarkPos, $bookmarkPos))
